$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# --- Crime statistics table updates (rows 15-29) ---
# Row 15
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -83.333333333333
$ws.Range("N15").Value = -83.333333333333

# Row 16
$ws.Range("D16").Value = 4
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -75
$ws.Range("J16").Value = 21
$ws.Range("K16").Value = -28.571428571428
$ws.Range("M16").Value = -51.612903225806
$ws.Range("N16").Value = -87.603305785124

# Row 17
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -53.846153846153
$ws.Range("I17").Value = 37
$ws.Range("J17").Value = 46
$ws.Range("K17").Value = -19.565217391304
$ws.Range("L17").Value = -13.953488372093
$ws.Range("M17").Value = -5.128205128205
$ws.Range("N17").Value = -57.471264367816

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("H18").Value = -33.333333333333
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = -41.666666666666
$ws.Range("L18").Value = -30
$ws.Range("N18").Value = -93.577981651376

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 300
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 10
$ws.Range("I19").Value = 71
$ws.Range("J19").Value = 57
$ws.Range("K19").Value = 24.561403508771
$ws.Range("L19").Value = 39.215686274509
$ws.Range("M19").Value = 42
$ws.Range("N19").Value = 5.970149253731

# Row 20
$ws.Range("C20").Value = "0"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = -26.315789473684
$ws.Range("L20").Value = 55.555555555555
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -90.54054054054

# Row 21
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = -61.538461538461
$ws.Range("F21").Value = 24
$ws.Range("G21").Value = 41
$ws.Range("H21").Value = -41.463414634146
$ws.Range("I21").Value = 152
$ws.Range("J21").Value = 173
$ws.Range("K21").Value = -12.138728323699
$ws.Range("L21").Value = 10.144927536231
$ws.Range("M21").Value = -1.298701298701
$ws.Range("N21").Value = -76.651305683563

# Row 22
$ws.Range("C22").Value = "0"

# Row 23
$ws.Range("D23").Value = "0"
$ws.Range("E23").Value = "***.*"
$ws.Range("G23").Value = 2
$ws.Range("L23").Value = -36.363636363636
$ws.Range("M23").Value = -53.333333333333

# Row 24
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 29
$ws.Range("G24").Value = 29
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 200
$ws.Range("J24").Value = 209
$ws.Range("K24").Value = -4.306220095693
$ws.Range("L24").Value = 20.481927710843
$ws.Range("M24").Value = 85.185185185185

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 36.842105263157
$ws.Range("I25").Value = 98
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 22.5
$ws.Range("L25").Value = 28.947368421052
$ws.Range("M25").Value = -40.60606060606

# Row 26
$ws.Range("D26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = -66.666666666666
$ws.Range("L26").Value = -66.666666666666

# Row 27
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60

# Row 28
$ws.Range("L28").Value = -33.333333333333

# Row 29
$ws.Range("L29").Value = -66.666666666666
